$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.088803768157959
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 3.853962182998657
$ws.Range("D1").Value = 3.295597791671753
$ws.Range("E1").Value = 1.687453746795654
